{"js": "const body = context.document.body;\n\n// Helper: locate the single range in the body whose text exactly matches\n// `oldText` and swap it for `newText`, via range.insertText(newText,\n// \"Replace\"). This rewrites only the matched range's text and keeps the\n// surrounding run/paragraph formatting (bold, italic, breaks, styles, etc.)\n// untouched.\nasync function replaceParagraphText(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error('Expected exactly 1 match for: ' + oldText.substring(0, 40) + ' but found ' + results.items.length);\n  }\n  results.items[0].insertText(newText, 'Replace');\n  await context.sync();\n}\n\n// 1) Ativacao date 2020 -> 2022\nawait replaceParagraphText(\"Ativa\u00e7\u00e3o: 01/01/2020\", \"Ativa\u00e7\u00e3o: 01/01/2022\");\n\n// 2) Remove 'Dupla Filtracao; ' from PT summary paragraph\nawait replaceParagraphText(\"Sistemas de Abastecimento de \u00c1gua; Tecnologias de Tratamento de \u00c1gua; Tratamento de \u00c1gua em Ciclo Completo; Desinfec\u00e7\u00e3o; Filtra\u00e7\u00e3o Direta Ascendente; Filtra\u00e7\u00e3o Direta Descendente; Dupla Filtra\u00e7\u00e3o; Floto-Filtra\u00e7\u00e3o; Filtra\u00e7\u00e3o em M\u00faltiplas Etapas; Tratamento dos Res\u00edduos Gerados nas Esta\u00e7\u00f5es de Tratamento de \u00c1gua.\", \"Sistemas de Abastecimento de \u00c1gua; Tecnologias de Tratamento de \u00c1gua; Tratamento de \u00c1gua em Ciclo Completo; Desinfec\u00e7\u00e3o; Filtra\u00e7\u00e3o Direta Ascendente; Filtra\u00e7\u00e3o Direta Descendente; Floto-Filtra\u00e7\u00e3o; Filtra\u00e7\u00e3o em M\u00faltiplas Etapas; Tratamento dos Res\u00edduos Gerados nas Esta\u00e7\u00f5es de Tratamento de \u00c1gua.\");\n\n// 3) Remove 'Double Filtration; ' from EN summary paragraph\nawait replaceParagraphText(\"Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Double Filtration; Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations.\", \"Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations.\");\n\n// 4) Replace full PT 'Programa' paragraph\nawait replaceParagraphText(\"dimensionamento de redes de distribui\u00e7\u00e3o de \u00e1gua; Mananciais Superficiais e Subterr\u00e2neos e Tecnologias de Tratamento de \u00c1gua; 1.1. Caracter\u00edsticas das \u00e1guas de interesse para o tratamento: caracter\u00edsticas f\u00edsicas, qu\u00edmicas e bacteriol\u00f3gicas; 1.2. Classifica\u00e7\u00e3o das \u00e1guas naturais destinadas ao abastecimento (CONAMA 20, NBR 12 216); 1.3. Padr\u00e3o de Potabilidade (Portaria 518/2004); 1.4. Tecnologias de Tratamento de \u00c1gua; 2. Unidades Constituintes de um Sistema de Abastecimento de \u00c1gua; 2.1. \u00c1guas Subterr\u00e2neas; 2.2. Capta\u00e7\u00e3o de \u00c1guas Superficiais - gradeamento, remo\u00e7\u00e3o de areia, casa de bombas; 2.3. Adu\u00e7\u00e3o (por gravidade, por recalque); 2.4. Reserva\u00e7\u00e3o; 2.5. Redes de distribui\u00e7\u00e3o: tipos de rede, pe\u00e7as e \u00f3rg\u00e3os acess\u00f3rios; 3. Sistema de Tratamento de \u00c1gua de Ciclo Completo; 3.1. Casa de Qu\u00edmica; 3.2. Coagula\u00e7\u00e3o e Mistura R\u00e1pida; caracter\u00edsticas de coagula\u00e7\u00e3o; tipos de coagulantes prim\u00e1rios e auxiliares de coagula\u00e7\u00e3o, flocula\u00e7\u00e3o e filtra\u00e7\u00e3o, diagramas de coagula\u00e7\u00e3o utilizando sais de alum\u00ednio e de ferro, tipos de unidades de mistura r\u00e1pida mecanizada e hidr\u00e1ulica, ensaios de Jarteste, projeto de unidade de mistura r\u00e1pida (misturados hidr\u00e1ulico; misturador mec\u00e2nico); 3.3. Flocula\u00e7\u00e3o: tipos de unidades de flocula\u00e7\u00e3o hidr\u00e1ulica e mecanizada, veicula\u00e7\u00e3o de \u00e1gua floculada nas ETAs, ensaios de Jarteste, projeto de unidades de flocula\u00e7\u00e3o (hidr\u00e1ulica; mec\u00e2nico); 3.4. Decanta\u00e7\u00e3o: decanta\u00e7\u00e3o convencional d e de alta taxa, projeto de dispositivos de entrada e sa\u00edda, sistema de remo\u00e7\u00e3o de lodo; 3.5. Flota\u00e7\u00e3o: conceitos e par\u00e2metros de projeto; dispositivos de remo\u00e7\u00e3o de lodo; 3.6. mecanismos da filtra\u00e7\u00e3o, materiais filtrantes e fundos de filtros, hidr\u00e1ulica da filtra\u00e7\u00e3o, filtra\u00e7\u00e3o com taxa constante e taxa declinante, modela\u00e7\u00e3o matem\u00e1tica para filtra\u00e7\u00e3o com taxa declinante, fluidifica\u00e7\u00e3o e expans\u00e3o de meio granulares, m\u00e9todos de lavagem de filtros, projetos de unidades de filtra\u00e7\u00e3o descendentes; 3.7. Desinfec\u00e7\u00e3o: principais desinfetantes, clora\u00e7\u00e3o e cloro-amonia\u00e7\u00e3o, pr\u00e9 e p\u00f3s-clora\u00e7\u00e3o, par\u00e2metros de projeto de c\u00e2mara de contato, subprodutos de desinfec\u00e7\u00e3o e principais desinfetantes alternativos, caracter\u00edsticas das cloraminas, oz\u00f4nio, per\u00f3xido de hidrog\u00eanio e di\u00f3xido de cloro; 3.8. Tratamento de res\u00edduos geradois nas ETAs e reuso de \u00e1gua recuperada: caracter\u00edsticas da \u00e1gua de lavagem dos filtros, descargas dos decantadores e de flotadores, clarifica\u00e7\u00e3o por sedimenta\u00e7\u00e3o, adensamento mec\u00e2nico, por gravidade e flota\u00e7\u00e3o, desaguamento por gravidade e mec\u00e2nico, propriedades do lodo e sua disposi\u00e7\u00e3o; 4. Tecnologias Alternativas de Tratamento de \u00c1gua; 4.1. Filtra\u00e7\u00e3o direta descendente: descri\u00e7\u00e3o da t\u00e9cnica e suas variantes, caracter\u00edsticas da coagula\u00e7\u00e3o, coagulantes prim\u00e1rios e auxiliares de coagula\u00e7\u00e3o, flocula\u00e7\u00e3o e filtra\u00e7\u00e3o, otimiza\u00e7\u00e3o do processo; 4.2. Filtra\u00e7\u00e3o direta ascendente: descri\u00e7\u00e3o da tecnologia e hist\u00f3rico sobre sua evolu\u00e7\u00e3o, caracter\u00edsticas da coagula\u00e7\u00e3o, variantes da tecnologia, m\u00e9todos de opera\u00e7\u00e3o: com e sem execu\u00e7\u00e3o de descargas de fundo intermedi\u00e1rias; 4.3. Dupla filtra\u00e7\u00e3o: caracter\u00edsticas principais da instala\u00e7\u00e3o, caracter\u00edsticas da coagula\u00e7\u00e3o, variantes da tecnologia e m\u00e9todos de opera\u00e7\u00e3o; 4.4. Filtra\u00e7\u00e3o em m\u00faltipas etapas - FiME: descri\u00e7\u00e3o geral da tecnologia, pr\u00e9-filtra\u00e7\u00e3o din\u00e2mica, pr\u00e9-filtra\u00e7\u00e3o em pedregulho com escoamento ascendente, descente ou horizontal, filtra\u00e7\u00e3o lenta em areia, filtra\u00e7\u00e3o lenta em areia e carv\u00e3o ativado granular; 4.5. Floto-filtra\u00e7\u00e3o: descri\u00e7\u00e3o geral da tecnologia, adequa\u00e7\u00e3o da filtra\u00e7\u00e3o r\u00e1pida ascendente com a flota\u00e7\u00e3o, remo\u00e7\u00e3o de lodo.\", \"- Tipos de dimensionamento de redes de distribui\u00e7\u00e3o de \u00e1gua;- Caracter\u00edsticas das \u00e1guas de interesse para o tratamento: caracter\u00edsticas f\u00edsicas, qu\u00edmicas e bacteriol\u00f3gicas; - Padr\u00e3o de Potabilidade; - Tecnologias de Tratamento de \u00c1gua;- Unidades Constituintes de um Sistema de Abastecimento de \u00c1gua;- Capta\u00e7\u00e3o de \u00c1guas Subterr\u00e2neas e Capta\u00e7\u00e3o de \u00c1guas Superficiais - Gradeamento, remo\u00e7\u00e3o de areia, casa de bombas; - Reserva\u00e7\u00e3o; - Redes de distribui\u00e7\u00e3o: tipos de rede, pe\u00e7as e \u00f3rg\u00e3os acess\u00f3rios; - Sistema de Tratamento de \u00c1gua de Ciclo Completo; - Coagula\u00e7\u00e3o-flocula\u00e7\u00e3o e Mistura R\u00e1pida; - Decanta\u00e7\u00e3o: decanta\u00e7\u00e3o convencional e de alta taxa e sistema de remo\u00e7\u00e3o de lodo;- Mecanismos da filtra\u00e7\u00e3o, materiais filtrantes e fundos de filtros, hidr\u00e1ulica da filtra\u00e7\u00e3o, filtra\u00e7\u00e3o com taxa constante e taxa declinante, - Desinfec\u00e7\u00e3o: principais desinfetantes, clora\u00e7\u00e3o e cloro-amonia\u00e7\u00e3o, pr\u00e9 e p\u00f3s-clora\u00e7\u00e3o, - Tratamento de res\u00edduos gerados nas ETAs e reuso de \u00e1gua recuperada\");\n\n// 5) Replace full EN 'Programa' paragraph (italic)\nawait replaceParagraphText(\"1. General View of a System of Water Supply, Use of sizing software for water distribution networks; Superficial and Subterranean Fountainheads and Water Treatment Technologies; 1.1. Water interest characteristics for treatment: physical, chemical and bacteriological characteristics; 1.2. Natural water classification destined to supply (CONAMA 20, NBR 12 216); 1.3. Potability Pattern (Portaria 518/2004); 1.4. Water Treatment Technologies; 2.System of Water Supply Constituent Units; 2.1. Subterranean Water; 2.2. Superficial Water Captivation - grating, sand removal, pump houses; 2.3. Delivery (through gravity, through suppression); 2.4. Reservation; 2.5. Distribution Systems: kinds of system, pieces and accessory organs; 3. System of Treatment Water Complete Cycle; 3.1. Chemistry House; 3.2. Coagulation and Fast Mixture; coagulation characteristics; types of primary coagulation and coagulation auxiliaries, flocculation and filtration, coagulation diagrams using aluminium and iron salt, types of mechanized and hydraulic fast mixture units,Jartest test tubes, fast mixture unit project (hydraulic mixer, mechanical mixer); 3.3. Flocculation: types of hydraulic and mechanized flocculation units, distribution of floccued water in the ETAs, Jarteste test tubes, flocculation unit project (hydraulic; mechanical); 3.4. Decantation: conventional and high tax decantation, input and output device project, sludge removal system; 3.5. Flotation: project concepts and parameters, sludge removal device; 3.6. filtration mechanisms, filtering materials and under gravel filters, filtration hydraulic, Constant and declining infiltration tax, fluidification and expansion of granular environments, filter washing methods, descedant filtration unit projects; 3.7. Disinfection: main disinfectants, chlorination and chlorine-ammoniation, pre and e pos- chlorination, Contact chamber parameters project, disinfection subprodcuts and main alternative disinfectants, chloramines, ozone, hydrogen peroxide and chlorine dioxide characteristics; 3.8. Treatment of waste generated in the ETAs and reuse of the recovered water: filter washing water characteristics, decanters and floating filters discharge, clarification through sedimentation, mechanical densification, through gravity and flotation, drainage through gravity and mechanical, sludge properties and its classification; 4. Water Treatment Alternative Technologies; 4.1. Direct descendant filtration: description of the tecnique and its variants, coagulation characteristics, primary coagulants and coagulation auxiliares, flocculation and filtration, optimization process; 4.2. Direct ascendant filtration: description of the technology and historical about its evolution, coagulation characteristics, technology variants, operation methods: with and without the execution of intermediary ground discharges; 4.3. Double filtration: main characteristics of the installation, coagulation characteristics, tehcnology variants and operation methods; 4.4. Filtration in multiple levels - FiML: general description of the technology, pre dynamic filtration, pre filtration in gravel stone with ascendant, descendant or horizontal drainage, slow filtration in sand, slow filtration in sand and granular activated carbon; 4.5. Floto-filtration: general description of the technology, adjustment of the fast ascendant filtration with the flotation, sludge removal.\", \"- Types of dimensioning of water distribution networks;- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water\");\n\n// 6) Replace 'Metodo' avaliacao sentence\nawait replaceParagraphText(\"Aulas te\u00f3ricas e pr\u00e1ticas, utiliza\u00e7\u00e3o de softwares, trabalhos de campo e exerc\u00edcios dirigidos.Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.\", \"Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios, trabalhos pr\u00e1ticos e relat\u00f3rios.\");\n", "ps1": "# Applies the LOB1257.docx content edits via Word COM Find/Replace.\n$wdFindContinue = 1\n$wdReplaceAll   = 2\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $ok = $range.Find.Execute(\n        $findText,\n        $false, $false, $false, $false, $false,\n        $true, $wdFindContinue, $false,\n        $replaceText, $wdReplaceAll\n    )\n    if (-not $ok) {\n        throw \"Find/Replace did not match: $($findText.Substring(0, [Math]::Min(40, $findText.Length)))\"\n    }\n}\n\n# Ativacao date: 2020 -> 2022\nReplace-ExactText \"Ativa\u00e7\u00e3o: 01/01/2020\" \"Ativa\u00e7\u00e3o: 01/01/2022\"\n\n# Remove 'Dupla Filtracao; ' from PT summary paragraph\nReplace-ExactText \"Sistemas de Abastecimento de \u00c1gua; Tecnologias de Tratamento de \u00c1gua; Tratamento de \u00c1gua em Ciclo Completo; Desinfec\u00e7\u00e3o; Filtra\u00e7\u00e3o Direta Ascendente; Filtra\u00e7\u00e3o Direta Descendente; Dupla Filtra\u00e7\u00e3o; Floto-Filtra\u00e7\u00e3o; Filtra\u00e7\u00e3o em M\u00faltiplas Etapas; Tratamento dos Res\u00edduos Gerados nas Esta\u00e7\u00f5es de Tratamento de \u00c1gua.\" \"Sistemas de Abastecimento de \u00c1gua; Tecnologias de Tratamento de \u00c1gua; Tratamento de \u00c1gua em Ciclo Completo; Desinfec\u00e7\u00e3o; Filtra\u00e7\u00e3o Direta Ascendente; Filtra\u00e7\u00e3o Direta Descendente; Floto-Filtra\u00e7\u00e3o; Filtra\u00e7\u00e3o em M\u00faltiplas Etapas; Tratamento dos Res\u00edduos Gerados nas Esta\u00e7\u00f5es de Tratamento de \u00c1gua.\"\n\n# Remove 'Double Filtration; ' from EN summary paragraph\nReplace-ExactText \"Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Double Filtration; Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations.\" \"Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations.\"\n\n# Replace full PT 'Programa' paragraph body\nReplace-ExactText \"dimensionamento de redes de distribui\u00e7\u00e3o de \u00e1gua; Mananciais Superficiais e Subterr\u00e2neos e Tecnologias de Tratamento de \u00c1gua; 1.1. Caracter\u00edsticas das \u00e1guas de interesse para o tratamento: caracter\u00edsticas f\u00edsicas, qu\u00edmicas e bacteriol\u00f3gicas; 1.2. Classifica\u00e7\u00e3o das \u00e1guas naturais destinadas ao abastecimento (CONAMA 20, NBR 12 216); 1.3. Padr\u00e3o de Potabilidade (Portaria 518/2004); 1.4. Tecnologias de Tratamento de \u00c1gua; 2. Unidades Constituintes de um Sistema de Abastecimento de \u00c1gua; 2.1. \u00c1guas Subterr\u00e2neas; 2.2. Capta\u00e7\u00e3o de \u00c1guas Superficiais - gradeamento, remo\u00e7\u00e3o de areia, casa de bombas; 2.3. Adu\u00e7\u00e3o (por gravidade, por recalque); 2.4. Reserva\u00e7\u00e3o; 2.5. Redes de distribui\u00e7\u00e3o: tipos de rede, pe\u00e7as e \u00f3rg\u00e3os acess\u00f3rios; 3. Sistema de Tratamento de \u00c1gua de Ciclo Completo; 3.1. Casa de Qu\u00edmica; 3.2. Coagula\u00e7\u00e3o e Mistura R\u00e1pida; caracter\u00edsticas de coagula\u00e7\u00e3o; tipos de coagulantes prim\u00e1rios e auxiliares de coagula\u00e7\u00e3o, flocula\u00e7\u00e3o e filtra\u00e7\u00e3o, diagramas de coagula\u00e7\u00e3o utilizando sais de alum\u00ednio e de ferro, tipos de unidades de mistura r\u00e1pida mecanizada e hidr\u00e1ulica, ensaios de Jarteste, projeto de unidade de mistura r\u00e1pida (misturados hidr\u00e1ulico; misturador mec\u00e2nico); 3.3. Flocula\u00e7\u00e3o: tipos de unidades de flocula\u00e7\u00e3o hidr\u00e1ulica e mecanizada, veicula\u00e7\u00e3o de \u00e1gua floculada nas ETAs, ensaios de Jarteste, projeto de unidades de flocula\u00e7\u00e3o (hidr\u00e1ulica; mec\u00e2nico); 3.4. Decanta\u00e7\u00e3o: decanta\u00e7\u00e3o convencional d e de alta taxa, projeto de dispositivos de entrada e sa\u00edda, sistema de remo\u00e7\u00e3o de lodo; 3.5. Flota\u00e7\u00e3o: conceitos e par\u00e2metros de projeto; dispositivos de remo\u00e7\u00e3o de lodo; 3.6. mecanismos da filtra\u00e7\u00e3o, materiais filtrantes e fundos de filtros, hidr\u00e1ulica da filtra\u00e7\u00e3o, filtra\u00e7\u00e3o com taxa constante e taxa declinante, modela\u00e7\u00e3o matem\u00e1tica para filtra\u00e7\u00e3o com taxa declinante, fluidifica\u00e7\u00e3o e expans\u00e3o de meio granulares, m\u00e9todos de lavagem de filtros, projetos de unidades de filtra\u00e7\u00e3o descendentes; 3.7. Desinfec\u00e7\u00e3o: principais desinfetantes, clora\u00e7\u00e3o e cloro-amonia\u00e7\u00e3o, pr\u00e9 e p\u00f3s-clora\u00e7\u00e3o, par\u00e2metros de projeto de c\u00e2mara de contato, subprodutos de desinfec\u00e7\u00e3o e principais desinfetantes alternativos, caracter\u00edsticas das cloraminas, oz\u00f4nio, per\u00f3xido de hidrog\u00eanio e di\u00f3xido de cloro; 3.8. Tratamento de res\u00edduos geradois nas ETAs e reuso de \u00e1gua recuperada: caracter\u00edsticas da \u00e1gua de lavagem dos filtros, descargas dos decantadores e de flotadores, clarifica\u00e7\u00e3o por sedimenta\u00e7\u00e3o, adensamento mec\u00e2nico, por gravidade e flota\u00e7\u00e3o, desaguamento por gravidade e mec\u00e2nico, propriedades do lodo e sua disposi\u00e7\u00e3o; 4. Tecnologias Alternativas de Tratamento de \u00c1gua; 4.1. Filtra\u00e7\u00e3o direta descendente: descri\u00e7\u00e3o da t\u00e9cnica e suas variantes, caracter\u00edsticas da coagula\u00e7\u00e3o, coagulantes prim\u00e1rios e auxiliares de coagula\u00e7\u00e3o, flocula\u00e7\u00e3o e filtra\u00e7\u00e3o, otimiza\u00e7\u00e3o do processo; 4.2. Filtra\u00e7\u00e3o direta ascendente: descri\u00e7\u00e3o da tecnologia e hist\u00f3rico sobre sua evolu\u00e7\u00e3o, caracter\u00edsticas da coagula\u00e7\u00e3o, variantes da tecnologia, m\u00e9todos de opera\u00e7\u00e3o: com e sem execu\u00e7\u00e3o de descargas de fundo intermedi\u00e1rias; 4.3. Dupla filtra\u00e7\u00e3o: caracter\u00edsticas principais da instala\u00e7\u00e3o, caracter\u00edsticas da coagula\u00e7\u00e3o, variantes da tecnologia e m\u00e9todos de opera\u00e7\u00e3o; 4.4. Filtra\u00e7\u00e3o em m\u00faltipas etapas - FiME: descri\u00e7\u00e3o geral da tecnologia, pr\u00e9-filtra\u00e7\u00e3o din\u00e2mica, pr\u00e9-filtra\u00e7\u00e3o em pedregulho com escoamento ascendente, descente ou horizontal, filtra\u00e7\u00e3o lenta em areia, filtra\u00e7\u00e3o lenta em areia e carv\u00e3o ativado granular; 4.5. Floto-filtra\u00e7\u00e3o: descri\u00e7\u00e3o geral da tecnologia, adequa\u00e7\u00e3o da filtra\u00e7\u00e3o r\u00e1pida ascendente com a flota\u00e7\u00e3o, remo\u00e7\u00e3o de lodo.\" \"- Tipos de dimensionamento de redes de distribui\u00e7\u00e3o de \u00e1gua;- Caracter\u00edsticas das \u00e1guas de interesse para o tratamento: caracter\u00edsticas f\u00edsicas, qu\u00edmicas e bacteriol\u00f3gicas; - Padr\u00e3o de Potabilidade; - Tecnologias de Tratamento de \u00c1gua;- Unidades Constituintes de um Sistema de Abastecimento de \u00c1gua;- Capta\u00e7\u00e3o de \u00c1guas Subterr\u00e2neas e Capta\u00e7\u00e3o de \u00c1guas Superficiais - Gradeamento, remo\u00e7\u00e3o de areia, casa de bombas; - Reserva\u00e7\u00e3o; - Redes de distribui\u00e7\u00e3o: tipos de rede, pe\u00e7as e \u00f3rg\u00e3os acess\u00f3rios; - Sistema de Tratamento de \u00c1gua de Ciclo Completo; - Coagula\u00e7\u00e3o-flocula\u00e7\u00e3o e Mistura R\u00e1pida; - Decanta\u00e7\u00e3o: decanta\u00e7\u00e3o convencional e de alta taxa e sistema de remo\u00e7\u00e3o de lodo;- Mecanismos da filtra\u00e7\u00e3o, materiais filtrantes e fundos de filtros, hidr\u00e1ulica da filtra\u00e7\u00e3o, filtra\u00e7\u00e3o com taxa constante e taxa declinante, - Desinfec\u00e7\u00e3o: principais desinfetantes, clora\u00e7\u00e3o e cloro-amonia\u00e7\u00e3o, pr\u00e9 e p\u00f3s-clora\u00e7\u00e3o, - Tratamento de res\u00edduos gerados nas ETAs e reuso de \u00e1gua recuperada\"\n\n# Replace full EN 'Programa' paragraph body (italic)\nReplace-ExactText \"1. General View of a System of Water Supply, Use of sizing software for water distribution networks; Superficial and Subterranean Fountainheads and Water Treatment Technologies; 1.1. Water interest characteristics for treatment: physical, chemical and bacteriological characteristics; 1.2. Natural water classification destined to supply (CONAMA 20, NBR 12 216); 1.3. Potability Pattern (Portaria 518/2004); 1.4. Water Treatment Technologies; 2.System of Water Supply Constituent Units; 2.1. Subterranean Water; 2.2. Superficial Water Captivation - grating, sand removal, pump houses; 2.3. Delivery (through gravity, through suppression); 2.4. Reservation; 2.5. Distribution Systems: kinds of system, pieces and accessory organs; 3. System of Treatment Water Complete Cycle; 3.1. Chemistry House; 3.2. Coagulation and Fast Mixture; coagulation characteristics; types of primary coagulation and coagulation auxiliaries, flocculation and filtration, coagulation diagrams using aluminium and iron salt, types of mechanized and hydraulic fast mixture units,Jartest test tubes, fast mixture unit project (hydraulic mixer, mechanical mixer); 3.3. Flocculation: types of hydraulic and mechanized flocculation units, distribution of floccued water in the ETAs, Jarteste test tubes, flocculation unit project (hydraulic; mechanical); 3.4. Decantation: conventional and high tax decantation, input and output device project, sludge removal system; 3.5. Flotation: project concepts and parameters, sludge removal device; 3.6. filtration mechanisms, filtering materials and under gravel filters, filtration hydraulic, Constant and declining infiltration tax, fluidification and expansion of granular environments, filter washing methods, descedant filtration unit projects; 3.7. Disinfection: main disinfectants, chlorination and chlorine-ammoniation, pre and e pos- chlorination, Contact chamber parameters project, disinfection subprodcuts and main alternative disinfectants, chloramines, ozone, hydrogen peroxide and chlorine dioxide characteristics; 3.8. Treatment of waste generated in the ETAs and reuse of the recovered water: filter washing water characteristics, decanters and floating filters discharge, clarification through sedimentation, mechanical densification, through gravity and flotation, drainage through gravity and mechanical, sludge properties and its classification; 4. Water Treatment Alternative Technologies; 4.1. Direct descendant filtration: description of the tecnique and its variants, coagulation characteristics, primary coagulants and coagulation auxiliares, flocculation and filtration, optimization process; 4.2. Direct ascendant filtration: description of the technology and historical about its evolution, coagulation characteristics, technology variants, operation methods: with and without the execution of intermediary ground discharges; 4.3. Double filtration: main characteristics of the installation, coagulation characteristics, tehcnology variants and operation methods; 4.4. Filtration in multiple levels - FiML: general description of the technology, pre dynamic filtration, pre filtration in gravel stone with ascendant, descendant or horizontal drainage, slow filtration in sand, slow filtration in sand and granular activated carbon; 4.5. Floto-filtration: general description of the technology, adjustment of the fast ascendant filtration with the flotation, sludge removal.\" \"- Types of dimensioning of water distribution networks;- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water\"\n\n# Replace avaliacao 'Metodo' sentence\nReplace-ExactText \"Aulas te\u00f3ricas e pr\u00e1ticas, utiliza\u00e7\u00e3o de softwares, trabalhos de campo e exerc\u00edcios dirigidos.Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.\" \"Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios, trabalhos pr\u00e1ticos e relat\u00f3rios.\"\n"}
